$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# --- Header row tweaks ---------------------------------------------------
# "horas extras" column header becomes "horas laboradas"
$ws.Range("D2").Value = "horas laboradas"

# --- Row 3: formula simplification + new "estado" value ------------------
$ws.Range("D3").Formula = "=C3-B3"

# E3 ("estado" column) gets a new value "TARDANZA"; reuse the time style
# (same numeric format as O2/P2/Q2, numFmtId 20) so no new style gets
# created in the workbook - just like copying formatting across in Excel.
$ws.Range("O2").Copy() | Out-Null
$ws.Range("E3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E3").Value = "TARDANZA"

# --- New attendance rows 4-8 ----------------------------------------------
# Each row: fecha (A), hora ingreso (B), hora salida (C), horas laboradas
# formula (D) and estado (E). Formats for A-D are copied from row 3 so the
# new rows reuse the existing date/time styles instead of creating new
# ones. Row 4's "estado" cell keeps the TARDANZA style (s=2); rows 5-8 use
# the worksheet's default (unstyled) formatting, matching the source data.
$rows = @(
    @{ Row = 4; Date = 45860; In = 0.3125;              Out = 0.75;                 Estado = "ASISTIO";  EstadoStyled = $true },
    @{ Row = 5; Date = 45861; In = 0.32361111111111113; Out = 0.71527777777777779;  Estado = "TARDANZA"; EstadoStyled = $false },
    @{ Row = 6; Date = 45862; In = 0.30902777777777779; Out = 0.70138888888888884;  Estado = "ASISTIO";  EstadoStyled = $false },
    @{ Row = 7; Date = 45863; In = 0.33333333333333331; Out = 0.75;                 Estado = "TARDANZA"; EstadoStyled = $false },
    @{ Row = 8; Date = 45773; In = 0.30902777777777779; Out = 0.63958333333333328;  Estado = "ASISTIO";  EstadoStyled = $false }
)

foreach ($r in $rows) {
    $n = $r.Row

    $ws.Range("A3").Copy() | Out-Null
    $ws.Range("A$n").PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range("A$n").Value = $r.Date

    $ws.Range("B3").Copy() | Out-Null
    $ws.Range("B$n").PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range("B$n").Value = $r.In

    $ws.Range("C3").Copy() | Out-Null
    $ws.Range("C$n").PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range("C$n").Value = $r.Out

    $ws.Range("D3").Copy() | Out-Null
    $ws.Range("D$n").PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range("D$n").Formula = "=C$n-B$n"

    if ($r.EstadoStyled) {
        $ws.Range("E3").Copy() | Out-Null
        $ws.Range("E$n").PasteSpecial($xlPasteFormats) | Out-Null
        $ws.Range("E$n").Value = $r.Estado
    } else {
        $ws.Range("E$n").Value = $r.Estado
    }
}

# --- Column width adjustments ---------------------------------------------
# The runtime's ColumnWidth setter adds a constant 5/6-character padding
# before storing the OOXML <col> width, so subtract it here to land on the
# target widths (B->11.6640625, C->10.21875, D->16.77734375, E->11).
$ws.Columns.Item(2).ColumnWidth = 10.830729166666666
$ws.Columns.Item(3).ColumnWidth = 9.385416666666666
$ws.Columns.Item(4).ColumnWidth = 15.944010416666666
$ws.Columns.Item(5).ColumnWidth = 10.166666666666666

# --- Selection --------------------------------------------------------------
$ws.Range("D2:D8").Select() | Out-Null
